# "Add more tests & include testing setups"
#
# - Update the multi-threaded benchmark numbers for a few plans and add a
#   new benchmark row (EU PebbleHost Extreme (9900K) 6GB) to the testing
#   data on Sheet1.
# - Remove the redundant/untitled duplicate chart that lived on Sheet1
#   (it duplicated the "Single-Threaded Performance" chart on Sheet2 but
#   had no title of its own).
# - Leave the UI selection where the author left it on each sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Update existing benchmark figures (Sheet1, "MULTI" table, col F) ---
$ws1.Range("F13").Value = 165
$ws1.Range("F14").Value = 182
$ws1.Range("F16").Value = 123
$ws1.Range("F17").Value = 124

# --- Add the new test row (19) to the MULTI table ---
$ws1.Range("E19").Value = "EU PebbleHost Extreme (9900K) 6GB"
$ws1.Range("F19").Value = 95
$ws1.Range("G19").Value = 0

# --- Remove the stray untitled chart embedded on Sheet1 (duplicate of the
#     "Single-Threaded Performance" chart, left over without its own title) ---
$charts1 = $ws1.ChartObjects()
for ($i = $charts1.Count; $i -ge 1; $i--) {
    $co = $charts1.Item($i)
    $title = ""
    if ($co.Chart.HasTitle) {
        $title = $co.Chart.ChartTitle.Text
    }
    if ($title.Trim() -eq "") {
        $co.Delete()
    }
}

# --- Restore each sheet's UI selection state ---
$ws1.Range("J16").Select()
$ws2.Range("W20").Select()
